$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.537.51'
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.377.99'
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.37'
$ws.Range("E5").Value = '  +5.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '187.06'
$ws.Range("E6").Value = '  -1.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.182'
$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.587'
$ws.Range("E10").Value = '  +0.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.50'
$ws.Range("E11").Value = '  +1.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000274'
$ws.Range("E12").Value = '  +1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.922.06'
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '637.12'
$ws.Range("E14").Value = '  +6.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.60'
$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.619.24'
$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.119'
$ws.Range("E17").Value = '  +0.95%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.376.97'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.02'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.17'
$ws.Range("E20").Value = '  +0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.910'
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.06'
$ws.Range("E22").Value = '  -3.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.11'
$ws.Range("E23").Value = '  +1.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.25'
$ws.Range("E24").Value = '  -0.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.02'
$ws.Range("E25").Value = '  +0.93%  '

$ws.Range("E26").Value = '  +3.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.54'
$ws.Range("E28").Value = '  +6.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.68'
$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("E30").Value = '  +2.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '610.28'
$ws.Range("E31").Value = '  +5.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.85'
$ws.Range("E32").Value = '  -2.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.012.49'
$ws.Range("E33").Value = '  +8.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.10'
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("E36").Value = '  +0.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.21'
$ws.Range("E37").Value = '  -0.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.84'
$ws.Range("E38").Value = '  +5.96%  '

$ws.Range("E39").Value = '  +3.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.68'
$ws.Range("E40").Value = '  -0.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.25'
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0{0}0703' -f [char]0x2083
$ws.Range("E42").Value = '  -0.98%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.344'
$ws.Range("E43").Value = '  +0.22%  '

$ws.Range("E44").Value = '  +0.25%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0423'
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.129'
$ws.Range("E46").Value = '  +0.12%  '

$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("E48").Value = '  +0.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.37'
$ws.Range("E49").Value = '  +10.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("E50").Value = '  -17.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '128.19'
$ws.Range("E51").Value = '  +3.14%  '
